$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 - Bug ID 12
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = "mastermind combat triggered even after correct puzzle"
$ws.Cells.Item(13, 3).Value = "ricky"
$ws.Cells.Item(13, 5).Value = "puzzle logic removed from combat logic"
$ws.Cells.Item(13, 4).Value = "incorrect puzzle logic in combat logic"
$ws.Cells.Item(13, 6).Value = "fixed"

# Row 14 - Bug ID 13
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = "combat infinite loop"
$ws.Cells.Item(14, 3).Value = "ricky"
$ws.Cells.Item(14, 4).Value = "combat loop incorrectly structured"
$ws.Cells.Item(14, 5).Value = "refactored combat loop"
$ws.Cells.Item(14, 6).Value = "fixed"

# Row 15 - Bug ID 14
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = "can still trade with enemies after defeating"
$ws.Cells.Item(15, 3).Value = "ricky"
$ws.Cells.Item(15, 4).Value = "no enemy death check for trade logic"

$ws.Range("F15").Select()
